$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp string (A1): refreshed at 17:52 instead of 17:22
$ws.Cells.Item(1,1).Value = "Datos actualizados a 22 de Abril de 2020 a las 17:52"

# Row 4: Estados Unidos (updated counts)
$ws.Cells.Item(4,2).Value = 820600
$ws.Cells.Item(4,3).Value = 1856
$ws.Cells.Item(4,4).Value = 83203
$ws.Cells.Item(4,5).Value = 691430
$ws.Cells.Item(4,6).Value = 14016
$ws.Cells.Item(4,7).Value = 649
$ws.Cells.Item(4,8).Value = 45967

# Row 31: was Singapur, now Polonia (re-sorted by Casos totales)
$ws.Cells.Item(31,1).Value = "Polonia"
$ws.Cells.Item(31,2).Value = 10169
$ws.Cells.Item(31,3).Value = 313
$ws.Cells.Item(31,4).Value = 1513
$ws.Cells.Item(31,5).Value = 8230
$ws.Cells.Item(31,6).Value = 160
$ws.Cells.Item(31,7).Value = 25
$ws.Cells.Item(31,8).Value = 426

# Row 32: was Pakistan, now Singapur (re-sorted by Casos totales)
$ws.Cells.Item(32,1).Value = "Singapur"
$ws.Cells.Item(32,2).Value = 10141
$ws.Cells.Item(32,3).Value = 1016
$ws.Cells.Item(32,4).Value = 839
$ws.Cells.Item(32,5).Value = 9291
$ws.Cells.Item(32,6).Value = 27
$ws.Cells.Item(32,7).Value = 0
$ws.Cells.Item(32,8).Value = 11

# Row 33: was Polonia, now Pakistan (re-sorted by Casos totales)
$ws.Cells.Item(33,1).Value = "Pakistan"
$ws.Cells.Item(33,2).Value = 10076
$ws.Cells.Item(33,3).Value = 511
$ws.Cells.Item(33,4).Value = 2156
$ws.Cells.Item(33,5).Value = 7708
$ws.Cells.Item(33,6).Value = 58
$ws.Cells.Item(33,7).Value = 11
$ws.Cells.Item(33,8).Value = 212

# Row 43: Chequia (updated counts)
$ws.Cells.Item(43,2).Value = 7087
$ws.Cells.Item(43,3).Value = 54
$ws.Cells.Item(43,4).Value = 1989
$ws.Cells.Item(43,5).Value = 4890
$ws.Cells.Item(43,6).Value = 80
$ws.Cells.Item(43,7).Value = 7
$ws.Cells.Item(43,8).Value = 208

# Row 80: was Afganistan, now Cuba (re-sorted by Casos totales)
$ws.Cells.Item(80,1).Value = "Cuba"
$ws.Cells.Item(80,2).Value = 1189
$ws.Cells.Item(80,3).Value = 52
$ws.Cells.Item(80,4).Value = 341
$ws.Cells.Item(80,5).Value = 808
$ws.Cells.Item(80,6).Value = 18
$ws.Cells.Item(80,7).Value = 2
$ws.Cells.Item(80,8).Value = 40

# Row 81: was Camerun, now Afganistan (re-sorted by Casos totales)
$ws.Cells.Item(81,1).Value = "Afganistan"
$ws.Cells.Item(81,2).Value = 1176
$ws.Cells.Item(81,3).Value = 84
$ws.Cells.Item(81,4).Value = 166
$ws.Cells.Item(81,5).Value = 970
$ws.Cells.Item(81,6).Value = 7
$ws.Cells.Item(81,7).Value = 4
$ws.Cells.Item(81,8).Value = 40

# Row 82: was Ghana, now Camerun (re-sorted by Casos totales)
$ws.Cells.Item(82,1).Value = "Camerun"
$ws.Cells.Item(82,2).Value = 1163
$ws.Cells.Item(82,3).Value = 0
$ws.Cells.Item(82,4).Value = 331
$ws.Cells.Item(82,5).Value = 789
$ws.Cells.Item(82,6).Value = 33
$ws.Cells.Item(82,7).Value = 0
$ws.Cells.Item(82,8).Value = 43

# Row 83: was Cuba, now Ghana (re-sorted by Casos totales)
$ws.Cells.Item(83,1).Value = "Ghana"
$ws.Cells.Item(83,2).Value = 1154
$ws.Cells.Item(83,3).Value = 112
$ws.Cells.Item(83,4).Value = 99
$ws.Cells.Item(83,5).Value = 1046
$ws.Cells.Item(83,6).Value = 4
$ws.Cells.Item(83,7).Value = 0
$ws.Cells.Item(83,8).Value = 9

# Row 110: Georgia (updated counts)
$ws.Cells.Item(110,2).Value = 416
$ws.Cells.Item(110,3).Value = 8
$ws.Cells.Item(110,4).Value = 107
$ws.Cells.Item(110,5).Value = 304
$ws.Cells.Item(110,6).Value = 6
$ws.Cells.Item(110,7).Value = 1
$ws.Cells.Item(110,8).Value = 5

# Row 113: Mauricio (updated counts)
$ws.Cells.Item(113,2).Value = 329
$ws.Cells.Item(113,3).Value = 1
$ws.Cells.Item(113,4).Value = 261
$ws.Cells.Item(113,5).Value = 59
$ws.Cells.Item(113,6).Value = 3
$ws.Cells.Item(113,7).Value = 0
$ws.Cells.Item(113,8).Value = 9

# Row 114: was Sri Lanka, now Mayotte (re-sorted by Casos totales)
$ws.Cells.Item(114,1).Value = "Mayotte"
$ws.Cells.Item(114,2).Value = 326
$ws.Cells.Item(114,3).Value = 15
$ws.Cells.Item(114,4).Value = 125
$ws.Cells.Item(114,5).Value = 197
$ws.Cells.Item(114,6).Value = 4
$ws.Cells.Item(114,7).Value = 0
$ws.Cells.Item(114,8).Value = 4

# Row 115: was Guatemala, now Sri Lanka (re-sorted by Casos totales)
$ws.Cells.Item(115,1).Value = "Sri Lanka"
$ws.Cells.Item(115,2).Value = 323
$ws.Cells.Item(115,3).Value = 13
$ws.Cells.Item(115,4).Value = 105
$ws.Cells.Item(115,5).Value = 211
$ws.Cells.Item(115,6).Value = 2
$ws.Cells.Item(115,7).Value = 0
$ws.Cells.Item(115,8).Value = 7

# Row 116: was Montenegro, now Guatemala (re-sorted by Casos totales)
$ws.Cells.Item(116,1).Value = "Guatemala"
$ws.Cells.Item(116,2).Value = 316
$ws.Cells.Item(116,3).Value = 22
$ws.Cells.Item(116,4).Value = 24
$ws.Cells.Item(116,5).Value = 284
$ws.Cells.Item(116,6).Value = 3
$ws.Cells.Item(116,7).Value = 1
$ws.Cells.Item(116,8).Value = 8

# Row 117: was Mayotte, now Montenegro (re-sorted by Casos totales)
$ws.Cells.Item(117,1).Value = "Montenegro"
$ws.Cells.Item(117,2).Value = 315
$ws.Cells.Item(117,3).Value = 2
$ws.Cells.Item(117,4).Value = 116
$ws.Cells.Item(117,5).Value = 194
$ws.Cells.Item(117,6).Value = 7
$ws.Cells.Item(117,7).Value = 0
$ws.Cells.Item(117,8).Value = 5
